$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in completion status for the web-interface task rows
$ws.Range("C111").Value = "完成"
$ws.Range("C112").Value = "未完成"
$ws.Range("C113").Value = "未完成"
$ws.Range("C114").Value = "完成"
$ws.Range("C115").Value = "未完成"

# Update the summary cell with the final write-up
$ws.Range("A116").Value = "总结：对于一些接口写的不够完善或者考虑不周的地方进行了修改。"

# Move the selection to the summary block, matching the saved view state
$ws.Range("A116:D118").Select()
